$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.395701
$ws.Range("H2").Value = 1.187103
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 46.06383033333333
$ws.Range("N2").Value = 138.191491
$ws.Range("O2").Value = 0.6011809800814781
$ws.Range("P2").Value = 0.601180980081478
$ws.Range("Q2").Value = 18.22750372673033
$ws.Range("R2").Value = 164.047533540573
$ws.Range("S2").Value = 0.6011809800814781
$ws.Range("T2").Value = 0.601180980081478

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.395701
$ws.Range("H3").Value = 1.187103
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 24.118612
$ws.Range("N3").Value = 72.35583600000001
$ws.Range("O3").Value = 0.3147730159528759
$ws.Range("P3").Value = 0.3147730159528759
$ws.Range("Q3").Value = 9.543758887012002
$ws.Range("R3").Value = 85.89382998310802
$ws.Range("S3").Value = 0.3147730159528759
$ws.Range("T3").Value = 0.3147730159528759

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.395701
$ws.Range("H4").Value = 1.187103
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 6.439792666666666
$ws.Range("N4").Value = 19.319378
$ws.Range("O4").Value = 0.0840460039656461
$ws.Range("P4").Value = 0.0840460039656461
$ws.Range("Q4").Value = 2.548232397992667
$ws.Range("R4").Value = 22.934091581934
$ws.Range("S4").Value = 0.0840460039656461
$ws.Range("T4").Value = 0.0840460039656461
